# "add test pport for swc2"
# Inserts a new Port row ("P_Test" / "SWC1, Delegation") right after SWC2's
# Task rows on the SWC_Composition sheet, pushing the SW_Composition1 /
# ToplevelComposition port tables below it down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SWC_Composition")

# Open up a blank row 12 (just under SWC2's existing rows) to host the new
# port entry; everything from the old row 12 onward shifts down by one.
$ws.Rows.Item(12).Insert()

# Clone the look of the row just below (a "Port / Provider / GeneralInterface"
# detail row) onto the freshly inserted row so borders/fonts match the table.
$ws.Range("B14:J14").Copy()
$ws.Range("B12:J12").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B12").RowHeight = 15

# Fill in the new port's data.
$ws.Range("D12").Value = "Port"
$ws.Range("E12").Value = "P_Test"
$ws.Range("F12").Value = "Provider"
$ws.Range("I12").Value = "GeneralInterface"
$ws.Range("J12").Value = "SWC1, Delegation"

# Matches the cursor position left behind in the authored workbook.
$ws.Range("J27").Select()
